$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 16 de Agosto de 2020 a las 00:45"

# Row 4
$ws.Range("B4").Value = 5526045
$ws.Range("C4").Value = 49779
$ws.Range("D4").Value = 2897539
$ws.Range("E4").Value = 2455969
$ws.Range("G4").Value = 1002
$ws.Range("H4").Value = 172537

# Row 5
$ws.Range("D5").Value = 2404272
$ws.Range("E5").Value = 805592

# Row 11
$ws.Range("B11").Value = 456689
$ws.Range("C11").Value = 11578
$ws.Range("D11").Value = 274420
$ws.Range("E11").Value = 167459
$ws.Range("G11").Value = 318
$ws.Range("H11").Value = 14810

# Row 23
$ws.Range("E23").Value = 101264
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 30409

# Row 32
$ws.Range("B32").Value = 96336
$ws.Range("C32").Value = 116
$ws.Range("D32").Value = 58835
$ws.Range("E32").Value = 32360
$ws.Range("G32").Value = 17
$ws.Range("H32").Value = 5141

# Row 50
$ws.Range("B50").Value = 53577
$ws.Range("C50").Value = 1360
$ws.Range("D50").Value = 38945
$ws.Range("E50").Value = 13547
$ws.Range("G50").Value = 12
$ws.Range("H50").Value = 1085

# Row 52
$ws.Range("B52").Value = 48770
$ws.Range("C52").Value = 325
$ws.Range("D52").Value = 36290
$ws.Range("E52").Value = 11506
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 974

# Row 60
$ws.Range("B60").Value = 37551
$ws.Range("C60").Value = 120
$ws.Range("G60").Value = 7
$ws.Range("H60").Value = 1370

# Row 75
$ws.Range("D75").Value = 16540
$ws.Range("E75").Value = 1528

# Row 86
$ws.Range("B86").Value = 9965
$ws.Range("C86").Value = 57
$ws.Range("E86").Value = 847

# Row 91
$ws.Range("B91").Value = 8588
$ws.Range("C91").Value = 39
$ws.Range("D91").Value = 7893
$ws.Range("E91").Value = 642

# Row 93
$ws.Range("B93").Value = 8343
$ws.Range("C93").Value = 83
$ws.Range("D93").Value = 7210
$ws.Range("E93").Value = 1083

# Row 130
$ws.Range("B130").Value = 2352
$ws.Range("C130").Value = 59
$ws.Range("D130").Value = 1631
$ws.Range("E130").Value = 713

# Row 133
$ws.Range("A133").Value = "Benin"
$ws.Range("B133").Value = 2063
$ws.Range("C133").Value = 49
$ws.Range("D133").Value = 1690
$ws.Range("E133").Value = 334
$ws.Range("H133").Value = 39

# Row 134
$ws.Range("A134").Value = "Tunez"
$ws.Range("B134").Value = 2023
$ws.Range("C134").Value = 120
$ws.Range("D134").Value = 1327
$ws.Range("E134").Value = 642
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 54

# Row 136
$ws.Range("B136").Value = 1954
$ws.Range("C136").Value = 7
$ws.Range("D136").Value = 1505
$ws.Range("E136").Value = 380

# Row 153
$ws.Range("B153").Value = 1130
$ws.Range("C153").Value = 6
$ws.Range("D153").Value = 836
$ws.Range("E153").Value = 267
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 27

# Row 165
$ws.Range("B165").Value = 674
$ws.Range("C165").Value = 25
$ws.Range("D165").Value = 310
$ws.Range("E165").Value = 342

# Row 167
$ws.Range("A167").Value = "Trinidad yTobago"
$ws.Range("B167").Value = 497
$ws.Range("C167").Value = 71
$ws.Range("D167").Value = 139
$ws.Range("E167").Value = 348
$ws.Range("H167").Value = 10

# Row 168
$ws.Range("A168").Value = "Taiwan"
$ws.Range("B168").Value = 482
$ws.Range("C168").Value = 1
$ws.Range("D168").Value = 450
$ws.Range("E168").Value = 25
$ws.Range("H168").Value = 7
